# Katalog kamerového vybavení - add a new station block and bump the date.
#
# Layout model: each "station" occupies two rows - a content row (station
# name / inspection note / camera / optics / accessories / controller note)
# followed by a thin spacer row. The controller legend (column I) stays
# pinned to the first two rows of the table and always shows "Kontroler 1"
# / "Kontroler 2"; the free-form "Kontroler: ..." note box lives in column
# D of the *last* station block.
#
# Adding a new station therefore:
#   - inserts two fresh rows at the bottom, carrying over the previous
#     station's formatting and its "Kontroler: ..." note,
#   - clears that note (and the matching C/D placeholder look) from the
#     row it used to occupy, turning it into a blank placeholder like the
#     neighbouring Camera/Optics cells,
#   - leaves the I-column legend text alone, since it still belongs to the
#     top of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Bump the worksheet date stamp.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Projekt: `nDatum: 07.07.2024"

# ---------------------------------------------------------------------
# 2) Clone the formatting of the existing station block (rows 4:5, still
#    untouched at this point) into the freshly inserted rows 6:7.
# ---------------------------------------------------------------------
$ws.Range("A4:I5").Copy()
$ws.Range("A6:I7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Populate row 6 with the content the old row 4 used to hold (the new
#    "last station" keeps the free-form controller note).
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Název stanice"
$ws.Range("B6").Value = "- popis inspekce`n"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "Kontroler: Kontroler 2  (FH-2050)`n"
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "`n"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""

# Row 7 is the spacer row under the new station block; it stays blank,
# same as row 5 used to be.

# ---------------------------------------------------------------------
# 5) Clear the controller note out of row 4 - C4/D4 become plain blank
#    placeholders, matching the look of the Camera/Optics placeholders
#    (E4/F4). The station name/description (A4/B4) and the "Kontroler 1"
#    legend (I4) are unchanged.
# ---------------------------------------------------------------------
$ws.Range("E4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = ""

$ws.Range("F4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "`n"

# ---------------------------------------------------------------------
# 6) Columns I of the new content/spacer rows (6:7) carry no legend text;
#    give them the plain "no border" look used elsewhere for such blanks
#    (same font/alignment as the other placeholders, but no cell border).
# ---------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Value = ""
$ws.Range("I6").Borders.LineStyle = -4142

$ws.Range("I6").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = ""

Write-Host "Edit complete"
